$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1 - copy formatting (bold + border) from the
# neighbouring "sum" header (G1) so it matches the other header cells, then
# set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Populate the Save column values for rows 2-15 (plain numbers, no special
# style - matching the other numeric data columns).
$saveValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 1
    10 = 1
    11 = 0
    12 = 0
    13 = 1
    14 = 1
    15 = 1
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
